$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tickers = @("IMX-USD", "TAO-USD", "GRT-USD", "PEPE-USD", "MNT-USD")

$startRow = 434
for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $tickers[$i]
}
